$d = $word.ActiveDocument

$replacements = @(
    @{old="77×97="; new="55×67="},
    @{old="13×70="; new="37×22="},
    @{old="36×80="; new="91×92="},
    @{old="76×70="; new="97×48="},
    @{old="23×64="; new="32×13="},
    @{old="88×63="; new="96×89="},
    @{old="42×85="; new="81×96="},
    @{old="73×75="; new="87×52="},
    @{old="70×59="; new="80×66="},
    @{old="97×41="; new="32×78="},
    @{old="32×39="; new="74×22="},
    @{old="19×24="; new="78×31="},
    @{old="92×63="; new="33×62="},
    @{old="15×34="; new="69×25="},
    @{old="67×53="; new="19×47="},
    @{old="35×55="; new="96×12="},
    @{old="22×76="; new="28×22="},
    @{old="81×69="; new="25×65="},
    @{old="73×85="; new="84×15="},
    @{old="14×63="; new="99×64="},
    @{old="37×79="; new="23×16="},
    @{old="23×25="; new="79×43="},
    @{old="38×23="; new="67×44="},
    @{old="90×68="; new="28×76="},
    @{old="93×36="; new="76×77="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
